{"js": "// The document contains four paragraphs (a title block repeated across the\n// guide) that read:\n//   \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u041f\u0435\u0440\u0441\u0435\u0443\u0441 \u0442\u043e\u043a\u043e\u043c 2018. \u0433\u043e\u0434\u0438\u043d\u0435 \u043f\u043e\u0441\u043c\u0430\u0442\u0440\u0430\u043c\u043e 30. \u043e\u043a\u0442\u043e\u0431\u0440\u0430 - 8.\n//    \u043d\u043e\u0432\u0435\u043c\u0431\u0440\u0430 \u0438 29. \u043d\u043e\u0432\u0435\u043c\u0431\u0440\u0430 - 8. \u0434\u0435\u0446\u0435\u043c\u0431\u0440\u0430\"\n// (the first occurrence also carries a leading formatting run with a\n// single space, plus a now-unused bookmark \"_Hlk514861060\" wrapped\n// around part of the text).\n//\n// Each of those paragraphs must become a single, plain run reading:\n//   \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 Taurus: 16. \u0434\u043e 25. \u0458\u0430\u043d\u0443\u0430\u0440\u0430\"\n// i.e. keep the (still-untranslated) constellation word \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435\", drop\n// \"\u041f\u0435\u0440\u0441\u0435\u0443\u0441\" (old constellation name) and the old observing dates, and\n// insert the new constellation placeholder \"Taurus\" with its new dates.\n\nconst searchText = \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435\";\nconst results = context.document.body.search(searchText, {\n  matchCase: true,\n  matchWholeWord: false,\n});\nresults.load(\"items\");\nawait context.sync();\n\nconst newText = \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 Taurus: 16. \u0434\u043e 25. \u0458\u0430\u043d\u0443\u0430\u0440\u0430\";\n\n// Collect the paragraphs that own each match \u2014 several hits can in\n// principle land in the same paragraph, so dedupe on uniqueLocalId.\nconst candidateParagraphs = [];\nfor (let i = 0; i < results.items.length; i++) {\n  const paragraph = results.items[i].paragraphs.getFirst();\n  paragraph.load(\"uniqueLocalId\");\n  candidateParagraphs.push(paragraph);\n}\nawait context.sync();\n\nconst seen = new Set();\nconst paragraphs = [];\nfor (const paragraph of candidateParagraphs) {\n  if (!seen.has(paragraph.uniqueLocalId)) {\n    seen.add(paragraph.uniqueLocalId);\n    paragraphs.push(paragraph);\n  }\n}\n\nfor (const paragraph of paragraphs) {\n  // Wipe every existing run (text + the leftover bookmark) in the\n  // paragraph, then drop in one clean run with the new wording. This\n  // mirrors the diff, which collapses the whole paragraph body down to a\n  // single <w:r><w:t>...</w:t></w:r> with no run formatting.\n  paragraph.clear();\n  paragraph.insertText(newText, Word.InsertLocation.start);\n}\n\nawait context.sync();\n", "ps1": "# The document repeats (four times) a title-block paragraph that reads:\n#   \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 \u041f\u0435\u0440\u0441\u0435\u0443\u0441 \u0442\u043e\u043a\u043e\u043c 2018. \u0433\u043e\u0434\u0438\u043d\u0435 \u043f\u043e\u0441\u043c\u0430\u0442\u0440\u0430\u043c\u043e 30. \u043e\u043a\u0442\u043e\u0431\u0440\u0430 - 8.\n#    \u043d\u043e\u0432\u0435\u043c\u0431\u0440\u0430 \u0438 29. \u043d\u043e\u0432\u0435\u043c\u0431\u0440\u0430 - 8. \u0434\u0435\u0446\u0435\u043c\u0431\u0440\u0430\"\n# The very first occurrence also still carries the long-unused bookmark\n# \"_Hlk514861060\" around part of its text.\n#\n# Each of those four paragraphs must collapse down to a single plain run:\n#   \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 Taurus: 16. \u0434\u043e 25. \u0458\u0430\u043d\u0443\u0430\u0440\u0430\"\n# keeping the (not-yet-translated) word \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435\", dropping the old\n# constellation name \"\u041f\u0435\u0440\u0441\u0435\u0443\u0441\" and the old date range, and writing in the\n# new constellation placeholder \"Taurus\" with its new observing dates.\n\n$d = $word.ActiveDocument\n\n# Drop the stale bookmark first (if present) so it doesn't linger as a\n# dangling <w:bookmarkStart/> once the text around it is replaced.\ntry {\n    $d.Bookmarks(\"_Hlk514861060\").Delete()\n} catch {\n    # Bookmark already gone / not present - nothing to do.\n}\n\n$newText = \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435 Taurus: 16. \u0434\u043e 25. \u0458\u0430\u043d\u0443\u0430\u0440\u0430\"\n\n# Gather every paragraph whose text contains the (capitalised, case-\n# sensitive) constellation heading \"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435\" - there are exactly four of\n# them in this guide.\n$targets = @()\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Contains(\"\u0421\u0430\u0437\u0432\u0435\u0436\u0452\u0435\")) {\n        $targets += $p\n    }\n}\n\nforeach ($p in $targets) {\n    $r = $p.Range\n    # Exclude the trailing paragraph mark from the range so we only wipe\n    # the run content, then delete it and insert the fresh, unformatted\n    # replacement text in its place.\n    $r.SetRange($r.Start, $r.End - 1)\n    $r.Delete()\n    $r.InsertAfter($newText)\n}\n"}
